# "Improved costcenter select in compensations"
#
# Changes applied:
#  1. Update a leader name: "Jonas Lind" -> "Sofia Nilsson" (C5)
#  2. Update costcenter / amount numbers in the compensation rows 2-5
#  3. Clear out the now-unused sample rows 6-9 (matching the blank template rows below)
#  4. Give the whole data table (A1:I21) a white fill so the costcenter /
#     compensation columns are easier to read against the thin grid border

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename leader
$ws.Range("C5").Value = "Sofia Nilsson"

# 2. Update numeric values (costcenter codes in column D, compensation amounts in column G)
$ws.Range("D2").Value = 400
$ws.Range("G2").Value = 120

$ws.Range("G3").Value = 300

$ws.Range("G4").Value = 90

$ws.Range("D5").Value = 4
$ws.Range("G5").Value = 100

# 3. Clear the sample rows that are no longer needed
$ws.Range("A6:I9").ClearContents()

# 4. Apply a white background fill across the whole table so it stands out
#    against the thin border grid used for the costcenter table
$ws.Range("A1:I21").Interior.Color = 16777215

Write-Host "Edit applied"
